# Generate Report for Handback
# Updates the localization-status report after a handback: statuses move
# from "Ready for handoff" to "Handed back: in sync with en-US", and the
# per-language sheets gain the Latest Target File / Latest Handback File /
# Latest Handback DateTime values (with a hyperlink on the target file).

$wb = $excel.ActiveWorkbook

$mdFile      = "170a976a-1720-47b4-8170-4c357abdc0e5.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/oltest/blob/cca40fbefc9cae053a996fb2d670a86bf69c3583/e2e/170a976a-1720-47b4-8170-4c357abdc0e5.md"
$handedBack  = "Handed back: in sync with en-US"

# ---- Overview sheet: both language-status columns move to "handed back" ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $handedBack
$zh.Range("C3").Value = $handedBack

$zh.Range("I2").Value = $mdFile
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, $null, $null, $mdFile) | Out-Null
$zh.Range("I3").Value = $mdFile
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl, $null, $null, $mdFile) | Out-Null

$zh.Range("J2").Value = "170a976a-1720-47b4-8170-4c357abdc0e5.f05706c3fed05b3e509b9b2cb9dc4f6afa05ae66.zh-cn.xlf"
$zh.Range("J3").Value = "170a976a-1720-47b4-8170-4c357abdc0e5.f05706c3fed05b3e509b9b2cb9dc4f6afa05ae66.zh-cn.xlf"

$zh.Range("K2").Value = "2016-08-13 23:27:43"
$zh.Range("K3").Value = "2016-08-13 23:27:43"

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $handedBack
$de.Range("C3").Value = $handedBack

$de.Range("I2").Value = $mdFile
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, $null, $null, $mdFile) | Out-Null
$de.Range("I3").Value = $mdFile
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl, $null, $null, $mdFile) | Out-Null

$de.Range("J2").Value = "170a976a-1720-47b4-8170-4c357abdc0e5.f05706c3fed05b3e509b9b2cb9dc4f6afa05ae66.de-de.xlf"
$de.Range("J3").Value = "170a976a-1720-47b4-8170-4c357abdc0e5.f05706c3fed05b3e509b9b2cb9dc4f6afa05ae66.de-de.xlf"

$de.Range("K2").Value = "2016-08-13 23:27:53"
$de.Range("K3").Value = "2016-08-13 23:27:53"

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
